$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Lattice Size, L (B1) from 7 to 8
$ws.Range("B1").Value = 8

# Update Temperature, T (B5) from 2.2999999999999998 to 2.1
$ws.Range("B5").Value = 2.1

# Update the active cell selection to E9
$ws.Range("E9").Select()
